$wb = $excel.ActiveWorkbook

# --- "Programs cost and coverage": insert a new "Delayed cord clamping" row ---
$ws = $wb.Worksheets.Item("Programs cost and coverage")

# Insert a new row above the old row 5 ("Family planning"); this pushes
# "Family planning" and everything below it down by one row.
$ws.Rows.Item(5).Insert()

$ws.Range("A5").Value = "Delayed cord clamping"
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0.95
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = "Constant (default)"

# Nudge the new row's formatting back onto the already-used style indices
# (71/72) that every other data row in this table uses, instead of the
# fresh derived styles a bare Insert() would otherwise allocate.
$ws.Range("B5:E5").Borders.LineStyle = 1

# Keep the sheet's recorded sort range in sync with the extra row.
$so = $ws.Sort
$so.SortFields.Clear()
$so.SortFields.Add($ws.Range("A2:A38"))
$so.SetRange($ws.Range("A2:D38"))
$so.Apply()

# --- "Programs target population": selection moved to B3 ---
$ws2 = $wb.Worksheets.Item("Programs target population")
$ws2.Range("B3").Select()

# Restore the first sheet as the active one (matches the saved workbook view).
$wb.Worksheets.Item(1).Activate()
